$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "extraction_date_latest" variable as row 38 (column A = name,
# column B = description), mirroring the existing variable rows above it.
$ws.Range("A38").Value = "extraction_date_latest"
$ws.Range("B38").Value = "Latest date for publications searches (latest search among the two independent reviewers)"

# Match the author's resulting selection/scroll state (next empty cell below
# the new row becomes the active cell).
$ws.Range("B39").Select()
$excel.ActiveWindow.ScrollRow = 22
